$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new columns -------------------------------------------------
# New column C: "fly_cost" (pushes fuel litres/cost/last activity right)
$ws.Columns.Item(3).EntireColumn.Insert()
# New column F: "day cost per all people" (pushes last activity right again)
$ws.Columns.Item(6).EntireColumn.Insert()

# --- Insert the new daily rows ---------------------------------------------------
# New row 3: 2021-06-17 (wx delay day) - shifts old 06-18 row (and below) down
$ws.Rows.Item(3).EntireRow.Insert()
# New row 7: 2021-06-20 (return to JAV to overnight)
$ws.Rows.Item(7).EntireRow.Insert()
# New row 8: 2021-06-21 (wx delay day)
$ws.Rows.Item(8).EntireRow.Insert()

# --- Insert the two new grand-total rows after the existing "total" row ----------
$ws.Rows.Item(13).EntireRow.Insert()
$ws.Rows.Item(14).EntireRow.Insert()

# --- Header row -------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "date (YYYY-MM-DD)"
$ws.Cells.Item(1,2).Value = "fly_time"
$ws.Cells.Item(1,3).Value = "fly_cost"
$ws.Cells.Item(1,4).Value = "fuel consumption litres"
$ws.Cells.Item(1,5).Value = "fuel consumption cost kDKK"
$ws.Cells.Item(1,6).Value = "day cost per all people"
$ws.Cells.Item(1,7).Value = "last activity"

# --- Daily rows ---------------------------------------------------------------
$rows = @(
  @("2021-06-16", 3.2, 42.6, 1274, 12.9, 7.5, "transit to Greenland."),
  @("2021-06-17", 4,   53.6, 0,    0,    7.5, "wx delay day"),
  @("2021-06-18", 2.8, 37.5, 2414, 24.5, 15,  "return to SFJ to overnight"),
  @("2021-06-19", 2.3, 30.8, 1381, 14,   15,  "return to SFJ to overnight"),
  @("2021-06-20", 2,   27.3, 1363, 13.8, 15,  "return to JAV to overnight"),
  @("2021-06-21", 4,   53.6, 0,    0,    15,  "wx delay day"),
  @("2021-06-22", 3.7, 50.1, 3213, 32.6, 15,  "to UAK to overnight"),
  @("2021-06-23", 1.8, 24,   718,  7.3,  15,  "transit"),
  @("2021-06-24", 3.4, 44.9, 3134, 31.8, 15,  "transit to SFJ"),
  @("2021-06-25", 3.5, 46.6, 1391, 14.1, 15,  "transit back to Canada")
)

$r = 2
foreach ($row in $rows) {
  # The date column must stay plain text (not get auto-parsed into a date
  # serial number): force a text number format for the assignment, then
  # clear the formatting again so the cell ends up back at the default
  # (unstyled) cell style, same as the rest of the data rows.
  $dateCell = $ws.Cells.Item($r,1)
  $dateCell.NumberFormat = "@"
  $dateCell.Value = $row[0]
  $dateCell.ClearFormats()

  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $ws.Cells.Item($r,5).Value = $row[4]
  $ws.Cells.Item($r,6).Value = $row[5]
  $ws.Cells.Item($r,7).Value = $row[6]
  $r = $r + 1
}

# --- Totals row (12) ------------------------------------------------------------
$ws.Cells.Item(12,1).Value = "total"
$ws.Cells.Item(12,2).Value = 30.7
$ws.Cells.Item(12,3).Value = 411
$ws.Cells.Item(12,4).Value = 14888
$ws.Cells.Item(12,5).Value = 151
$ws.Cells.Item(12,6).Value = 135

# --- Grand total rows (13, 14) ---------------------------------------------------
$ws.Cells.Item(13,1).Value = "grand total (MDKK)"
$ws.Cells.Item(13,2).Value = 0.697

$ws.Cells.Item(14,1).Value = "grand total incl. quarantine (MDKK)"
$ws.Cells.Item(14,2).Value = 0.772
